# Add two new columns "I0" (I) and "IF" (J) to the sheet, mirroring the
# existing header/style pattern used by the other columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers (row 1) ---
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the header formatting used by the other header cells (e.g. H1):
# bold font, thin border all around, centered horizontally, top vertically.
$headerSample = $ws.Range("H1")
foreach ($addr in @("I1", "J1")) {
    $cell = $ws.Range($addr)
    $cell.Font.Bold = $headerSample.Font.Bold
    $cell.HorizontalAlignment = $headerSample.HorizontalAlignment
    $cell.VerticalAlignment = $headerSample.VerticalAlignment
    $cell.Borders.LineStyle = $headerSample.Borders.LineStyle
}

# --- Data rows (2-10) ---
$values = @{
    2  = 8
    3  = 9
    4  = 7
    5  = 8
    6  = 6
    7  = 9
    8  = 6
    9  = 9
    10 = 9
}

foreach ($row in $values.Keys) {
    $v = $values[$row]
    $ws.Cells.Item($row, 9).Value = $v   # column I
    $ws.Cells.Item($row, 10).Value = $v  # column J
}
